$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "minAcro" -> "Misc Values" (row 5, columns B and D mirror each other)
$ws.Range("B5").Value = "Misc Values"
$ws.Range("D5").Value = "Misc Values"

# Row 16/17: reword D16, and highlight B16/D16/D17 bold + light-blue.
$ws.Range("D16").Value = "Send accel calibration values"

$ws.Range("B16").Font.Bold = $true
$ws.Range("B16").Font.Color = 15773696
$ws.Range("D16").Font.Bold = $true
$ws.Range("D16").Font.Color = 15773696
$ws.Range("D17").Font.Bold = $true
$ws.Range("D17").Font.Color = 15773696

# Row 18: new content added in B and D, also bold/light-blue.
$ws.Range("B18").Value = "Write mag calibration values"
$ws.Range("D18").Value = "Send mag cal values"

$ws.Range("B18").Font.Bold = $true
$ws.Range("B18").Font.Color = 15773696
$ws.Range("D18").Font.Bold = $true
$ws.Range("D18").Font.Color = 15773696

# New values added to previously empty cells further down the D column.
$ws.Range("D20").Value = "Send Pressure Altitude"
$ws.Range("D24").Value = "Send Command in Detent Discretes"
$ws.Range("D26").Value = "Send 100 Hz loop time"

# Matches the selection left behind in the saved workbook.
$ws.Range("I30").Select()
